# Applies:
#  1. Remove the stray _GoBack bookmark pair sitting alone in the empty
#     "teichoice" paragraph right after "Leopold Volkmer / D. J. Pajek 1881."
#  2. Remove the direct <w:u w:val="none"/> override on the teiadd-styled
#     run "to" (ſe je *to* vdalo ...)
#  3. Remove the direct <w:u w:val="none"/> override on the teiadd-styled
#     run "oni" (nebi bli *oni* od povanya ...)
#  4. Split the teisurplus-styled run "ga" (dabi tega*ga*) into "g" + "a",
#     re-inserting the _GoBack bookmark pair between the two new runs.

$d = $word.ActiveDocument

# --- 1. Remove the old _GoBack bookmark -----------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# --- 2. Strip the direct underline override off the "to" run --------------
$rng = $d.Content
$rng.Find.Execute("e je to vdalo", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$toStart = $rng.Start + 5          # length of "e je "
$toEnd = $toStart + 2              # length of "to"
$toRng = $d.Range($toStart, $toEnd)
if ($toRng.Text -ne "to") { throw "unexpected text at 'to' range: $($toRng.Text)" }
$toRng.Delete()
$d.Range($toStart, $toStart).InsertAfter("to") | Out-Null
$d.Range($toStart, $toStart + 2).Style = $d.Styles("teiadd")

# --- 3. Strip the direct underline override off the "oni" run -------------
$rng2 = $d.Content
$rng2.Find.Execute("nebi bli oni od povanya", $true, $false, $false, $false, `
                    $false, $true, 1, $false, "", 0) | Out-Null
$oniStart = $rng2.Start + 9        # length of "nebi bli "
$oniEnd = $oniStart + 3            # length of "oni"
$oniRng = $d.Range($oniStart, $oniEnd)
if ($oniRng.Text -ne "oni") { throw "unexpected text at 'oni' range: $($oniRng.Text)" }
$oniRng.Delete()
$d.Range($oniStart, $oniStart).InsertAfter("oni") | Out-Null
$d.Range($oniStart, $oniStart + 3).Style = $d.Styles("teiadd")

# --- 4. Split "ga" (teisurplus) into "g" + bookmark + "a" ------------------
$rng3 = $d.Content
$rng3.Find.Execute("dabi tegaga", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$gaStart = $rng3.Start + 9         # length of "dabi tega"
$gaEnd = $gaStart + 2              # length of "ga"
$gaRng = $d.Range($gaStart, $gaEnd)
if ($gaRng.Text -ne "ga") { throw "unexpected text at 'ga' range: $($gaRng.Text)" }

# Drop the trailing "a" so the run's text becomes "g"
$d.Range($gaStart + 1, $gaEnd).Delete()

# Re-insert the _GoBack bookmark right after "g"
$d.Bookmarks.Add("_GoBack", $d.Range($gaStart + 1, $gaStart + 1)) | Out-Null

# Insert a fresh "a" run (after the bookmark) using the same teisurplus style
$d.Range($gaStart + 1, $gaStart + 1).InsertAfter("a") | Out-Null
$d.Range($gaStart + 1, $gaStart + 2).Style = $d.Styles("teisurplus")

Write-Output "edits applied"
